# Update the "RESULT" schedule sheet with the revised subject assignments.
# (commit: "TODO init_new + all TODOs")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Week 1 (rows 2-3)
$ws.Range("A2").Value = "Phys1-1"
$ws.Range("B2").Value = "Math1-2"
$ws.Range("C2").Value = "Phys1-3"
$ws.Range("D2").Value = "Litr1-4"
$ws.Range("E2").Value = "Math1-5"
$ws.Range("A3").Value = "Math1-1"
$ws.Range("B3").Value = "Litr1-2"
$ws.Range("C3").Value = "Litr1-3"
$ws.Range("D3").Value = "Math1-4"
$ws.Range("E3").Value = "Litr1-5"

# Week 2 (rows 6-7)
$ws.Range("A6").Value = "Math2-1"
$ws.Range("B6").Value = "English2-2"
$ws.Range("C6").Value = "Phys2-3"
$ws.Range("D6").Value = "English2-4"
$ws.Range("E6").Value = "Math2-5"
$ws.Range("A7").Value = "English2-1"
$ws.Range("B7").Value = "Phys2-2"
$ws.Range("C7").Value = "Math2-3"
$ws.Range("D7").Value = "Phys2-4"
$ws.Range("E7").Value = "Phys2-5"

# Week 3 (rows 10-12)
$ws.Range("A10").Value = "Phys3-1"
$ws.Range("B10").Value = "Phys3-2"
$ws.Range("C10").Value = "Russian3-3"
$ws.Range("D10").Value = "Phys3-4"
$ws.Range("E10").Value = "English3-5"
$ws.Range("A11").Value = "Russian3-1"
$ws.Range("B11").Value = "Russian3-2"
$ws.Range("C11").Value = "Phys3-3"
$ws.Range("D11").Value = "Litra3-4"
$ws.Range("E11").Value = "Litra3-5"
$ws.Range("A12").Value = "English3-1"

# Week 4 (rows 15-16)
$ws.Range("A15").Value = "Math4-1"
$ws.Range("B15").Value = "Phys4-2"
$ws.Range("C15").Value = "Math4-3"
$ws.Range("D15").Value = "Math4-4"
$ws.Range("E15").Value = "Math4-5"
$ws.Range("A16").Value = "English4-1"
$ws.Range("B16").Value = "Math4-2"
$ws.Range("C16").Value = "Phys4-3"
$ws.Range("D16").Value = "Phys4-4"
$ws.Range("E16").Value = "English4-5"

# Week 5 (row 19)
$ws.Range("A19").Value = "Russian5-1"
$ws.Range("B19").Value = "Phys5-2"
$ws.Range("C19").Value = "Phys5-3"
$ws.Range("D19").Value = "Phys5-4"
$ws.Range("E19").Value = "Phys5-5"
